# Commit: "Adição dos exercicios de lista"
#
# 1) On the "78" worksheet (Tabela3: Nome/Input/Nome Buscado/Encontrado/Output),
#    rename the last column header from "Output" to "Saida" and fix the final
#    "found" output cell from the placeholder "{ACHEI}" to "Achei!".
# 2) Add a brand-new, empty worksheet named "80" after the existing sheets
#    (78, 79, 82) to hold the next list exercise, and make it the active sheet.

$wb = $excel.ActiveWorkbook

# --- 1) Update the "78" sheet / Tabela3 ---
$ws78 = $wb.Worksheets.Item("78")
$ws78.Range("E1").Value = "Saida"
$ws78.Range("E15").Value = "Achei!"
$ws78.Range("E16").Select()

# --- 2) Append a new, blank worksheet "80" after the last sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws80 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws80.Name = "80"
$ws80.Range("A1:J16").Select()
